$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$officesText = @'
0.4% CR/LFINF+CDM/H:1/Offices
4.8% CR/LFINF+CDM/H:2/Offices
22.8% CR+PC/LWAL+CDM/HBET:3-5/Offices
0.2% CR+PC/LWAL+CDM/HBET:6-/Offices
1.3% CR/LFINF+CDN/HBET:3-5/Offices
0.1% CR/LFINF+CDN/HBET:6-/Offices
6.0% CR+PC/LWAL+CDN/HBET:3-5/Offices
8.4% MCF/LWAL+CDL/H:1/Offices
12.5% MCF/LWAL+CDL/H:2/Offices
1.7% MCF/LWAL+CDL/HBET:3-5/Offices
8.8% MUR+CL/LWAL+CDN/H:1/Offices
12.9% MUR+CL/LWAL+CDN/H:2/Offices
8.4% MUR+ST/LWAL+CDN/H:1/Offices
11.7% MUR+ST/LWAL+CDN/H:2/Offices
0.0% CR/LFINF+CDN/H:1/Offices
0.0% CR/LFINF+CDN/H:2/Offices
0.0% CR+PC/LWAL+CDN/H:2/Offices
0.0% W/LFM+CDM/H:1/Offices
0.0% W/LFM+CDM/H:2/Offices
0.0% W/LFM+CDN/H:1/Offices
0.0% W/LFM+CDN/H:2/Offices
0.0% MUR+ADO/LWAL+CDN/H:2/Offices
0.0% MUR+ADO/LWAL+CDN/H:1/Offices
0.0% W/LFM+CDM/HBET:3-5/Offices
'@

$tradeText = @'
21.8% CR/LFINF+CDM/H:1/Trade
1.5% CR/LFINF+CDM/H:2/Trade
0.0% CR+PC/LWAL+CDM/HBET:3-5/Trade
0.0% CR+PC/LWAL+CDM/HBET:6-/Trade
0.0% CR/LFINF+CDN/HBET:3-5/Trade
0.0% CR/LFINF+CDN/HBET:6-/Trade
0.0% CR+PC/LWAL+CDN/HBET:3-5/Trade
19.0% MCF/LWAL+CDL/H:1/Trade
1.6% MCF/LWAL+CDL/H:2/Trade
0.0% MCF/LWAL+CDL/HBET:3-5/Trade
18.0% MUR+CL/LWAL+CDN/H:1/Trade
0.7% MUR+CL/LWAL+CDN/H:2/Trade
18.2% MUR+ST/LWAL+CDN/H:1/Trade
0.0% MUR+ST/LWAL+CDN/H:2/Trade
12.2% CR/LFINF+CDN/H:1/Trade
0.1% CR/LFINF+CDN/H:2/Trade
5.4% CR+PC/LWAL+CDN/H:2/Trade
0.5% W/LFM+CDM/H:1/Trade
0.0% W/LFM+CDM/H:2/Trade
1.0% W/LFM+CDN/H:1/Trade
0.0% W/LFM+CDN/H:2/Trade
0.0% MUR+ADO/LWAL+CDN/H:2/Trade
0.0% MUR+ADO/LWAL+CDN/H:1/Trade
0.0% W/LFM+CDM/HBET:3-5/Trade
'@

$hotelsText = @'
0.0% CR/LFINF+CDM/H:1/Hotels
 0.1% CR/LFINF+CDM/H:2/Hotels
 30.5% CR+PC/LWAL+CDM/HBET:3-5/Hotels
 0.2% CR+PC/LWAL+CDM/HBET:6-/Hotels
 4.9% CR/LFINF+CDN/HBET:3-5/Hotels
 0.1% CR/LFINF+CDN/HBET:6-/Hotels
 6.4% CR+PC/LWAL+CDN/HBET:3-5/Hotels
 13.4% MCF/LWAL+CDL/H:1/Hotels
 1.9% MCF/LWAL+CDL/H:2/Hotels
 5.4% MCF/LWAL+CDL/HBET:3-5/Hotels
 15.8% MUR+CL/LWAL+CDN/H:1/Hotels
 0.8% MUR+CL/LWAL+CDN/H:2/Hotels
 12.6% MUR+ST/LWAL+CDN/H:1/Hotels
 2.3% MUR+ST/LWAL+CDN/H:2/Hotels
 0.0% CR/LFINF+CDN/H:1/Hotels
 0.0% CR/LFINF+CDN/H:2/Hotels
 0.0% CR+PC/LWAL+CDN/H:2/Hotels
 0.3% W/LFM+CDM/H:1/Hotels
 0.3% W/LFM+CDM/H:2/Hotels
 0.6% W/LFM+CDN/H:1/Hotels
 0.6% W/LFM+CDN/H:2/Hotels
 0.0% MUR+ADO/LWAL+CDN/H:2/Hotels
 3.8% MUR+ADO/LWAL+CDN/H:1/Hotels
 0.0% W/LFM+CDM/HBET:3-5/Hotels
'@

$ws.Range("B2").Value = $officesText
$ws.Range("D2").Value = $hotelsText
$ws.Range("C2").Value = $tradeText

$ws.Range("B2:D2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 380

$ws.Columns.Item(2).ColumnWidth = 34.666666666666664
$ws.Columns.Item(3).ColumnWidth = 33.5
$ws.Columns.Item(4).ColumnWidth = 41.333333333333336

$ws.Range("B2").Select()
